# v0.1: Updated profiles for better results
#
# The "Sn96.5/Ag3/Cu0.5" (first) reflow profile's time column (C6:C9 on
# sheet "Hoja1") is updated with new, lower time values so the resulting
# ramp-rate/profile gives better results. The dependent ramp-rate formula
# in column E recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C6").Value = 140
$ws.Range("C7").Value = 230
$ws.Range("C8").Value = 270
$ws.Range("C9").Value = 370

# Leave the selection where the author last left off editing.
[void]$ws.Range("F17").Select()
